$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '57.080.11'
$ws.Range('E2').Value = '  -1.90%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.067.71'
$ws.Range('E3').Value = '  -2.19%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '520.70'
$ws.Range('E5').Value = '  -1.81%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '135.29'
$ws.Range('E6').Value = '  -5.26%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.066.82'
$ws.Range('E8').Value = '  -1.93%  '
$ws.Range('E9').Value = '  +5.27%  '
$ws.Range('E10').Value = '  +2.02%  '
$ws.Range('E11').Value = '  -3.07%  '
$ws.Range('E12').Value = '  +2.01%  '
$ws.Range('E13').Value = '  +1.01%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.596.24'
$ws.Range('E14').Value = '  -1.84%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '25.14'
$ws.Range('E15').Value = '  -2.38%  '
$ws.Range('E16').Value = '  -3.51%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '57.147.11'
$ws.Range('E17').Value = '  -1.85%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.061.29'
$ws.Range('E18').Value = '  -2.29%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.87'
$ws.Range('E19').Value = '  -4.22%  '
$ws.Range('E20').Value = '  -3.30%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.80'
$ws.Range('E21').Value = '  -2.42%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '348.61'
$ws.Range('E22').Value = '  +1.52%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.998'
$ws.Range('E23').Value = '  -0.23%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '69.05'
$ws.Range('E24').Value = '  +2.06%  '
$ws.Range('E25').Value = '  -3.41%  '
$ws.Range('E26').Value = '  -2.50%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.999'
$ws.Range('E27').Value = '  -0.19%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.0₃0855'
$ws.Range('E28').Value = '  -8.51%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.998'
$ws.Range('E29').Value = '  -0.02%  '
$ws.Range('E30').Value = '  -3.07%  '
$ws.Range('E31').Value = '  -2.03%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.78'
$ws.Range('E32').Value = '  -9.87%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '20.87'
$ws.Range('E33').Value = '  -1.46%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '158.74'
$ws.Range('E34').Value = '  +0.09%  '
$ws.Range('E35').Value = '  +0.31%  '
$ws.Range('E36').Value = '  -5.69%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.97'
$ws.Range('E37').Value = '  -3.66%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '25.26'
$ws.Range('E38').Value = '  -4.01%  '
$ws.Range('E39').Value = '  -2.07%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0654'
$ws.Range('E40').Value = '  -2.10%  '
$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.56'
$ws.Range('E41').Value = '  -6.07%  '
$ws.Range('B42').Value = 'Filecoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.02'
$ws.Range('E42').Value = '  +0.49%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.691'
$ws.Range('E43').Value = '  -0.95%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.405.26'
$ws.Range('E44').Value = '  +6.19%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '36.44'
$ws.Range('E45').Value = '  -0.71%  '
$ws.Range('E46').Value = '  +0.09%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.108.75'
$ws.Range('E47').Value = '  -1.97%  '
$ws.Range('E48').Value = '  -1.07%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '5.96'
$ws.Range('E49').Value = '  -2.56%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.930'
$ws.Range('E50').Value = '  -7.26%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '19.45'
$ws.Range('E51').Value = '  -6.14%  '
